# This workbook contains a single "forecast series" table in Sheet1
# (columns: y_0 date, y_0 year, y_0 forecast, y_1 year, y_1 forecast).
# The commit re-bases the forecast on an extra leading observation year
# (2007), which shifts every existing row down by one and refreshes the
# forecast values throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2, shifting rows 2-18 down to 3-19.
$ws.Rows.Item(2).Insert()

# Match formatting used by the rest of the table: column A carries a
# bordered/centered custom date-time number format; B:E carry no explicit
# style. The freshly inserted row inherits an unrelated bold style from
# row 1, so reset B2:E2 to the default look and apply column A's look to A2.
$ws.Range("B2:E2").ClearFormats()
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Borders.LineStyle = 1
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160

# Rewrite the full data block (rows 2-19) with the refreshed series.
# Row 2
$ws.Cells.Item(2, 1).Value = 39400
$ws.Cells.Item(2, 2).Value = 2007
$ws.Cells.Item(2, 3).Value = 1.144978573787081
$ws.Cells.Item(2, 4).Value = 2008
$ws.Cells.Item(2, 5).Value = 2.918510996763723
# Row 3
$ws.Cells.Item(3, 1).Value = 39765
$ws.Cells.Item(3, 2).Value = 2008
$ws.Cells.Item(3, 3).Value = 1.381024225294869
$ws.Cells.Item(3, 4).Value = 2009
$ws.Cells.Item(3, 5).Value = 4.264380119800992
# Row 4
$ws.Cells.Item(4, 1).Value = 40130
$ws.Cells.Item(4, 2).Value = 2009
$ws.Cells.Item(4, 3).Value = -0.8792832172735965
$ws.Cells.Item(4, 4).Value = 2010
$ws.Cells.Item(4, 5).Value = 3.109784435759599
# Row 5
$ws.Cells.Item(5, 1).Value = 40494
$ws.Cells.Item(5, 2).Value = 2010
$ws.Cells.Item(5, 3).Value = 0.9337833426867448
$ws.Cells.Item(5, 4).Value = 2011
$ws.Cells.Item(5, 5).Value = 2.730731696345146
# Row 6
$ws.Cells.Item(6, 1).Value = 40862
$ws.Cells.Item(6, 2).Value = 2011
$ws.Cells.Item(6, 3).Value = 2.791140000794279
$ws.Cells.Item(6, 4).Value = 2012
$ws.Cells.Item(6, 5).Value = 1.683857142130885
# Row 7
$ws.Cells.Item(7, 1).Value = 41228
$ws.Cells.Item(7, 2).Value = 2012
$ws.Cells.Item(7, 3).Value = 0.4451370000809973
$ws.Cells.Item(7, 4).Value = 2013
$ws.Cells.Item(7, 5).Value = -0.2532347529486723
# Row 8
$ws.Cells.Item(8, 1).Value = 41592
$ws.Cells.Item(8, 2).Value = 2013
$ws.Cells.Item(8, 3).Value = 0.2545814083968478
$ws.Cells.Item(8, 4).Value = 2014
$ws.Cells.Item(8, 5).Value = 2.553470871380537
# Row 9
$ws.Cells.Item(9, 1).Value = 41957
$ws.Cells.Item(9, 2).Value = 2014
$ws.Cells.Item(9, 3).Value = 1.297015177357297
$ws.Cells.Item(9, 4).Value = 2015
$ws.Cells.Item(9, 5).Value = -0.434146007584113
# Row 10
$ws.Cells.Item(10, 1).Value = 42321
$ws.Cells.Item(10, 2).Value = 2015
$ws.Cells.Item(10, 3).Value = 1.365576377841027
$ws.Cells.Item(10, 4).Value = 2016
$ws.Cells.Item(10, 5).Value = 2.383242923544526
# Row 11
$ws.Cells.Item(11, 1).Value = 42689
$ws.Cells.Item(11, 2).Value = 2016
$ws.Cells.Item(11, 3).Value = 2.204449574611278
$ws.Cells.Item(11, 4).Value = 2017
$ws.Cells.Item(11, 5).Value = 1.688977015142101
# Row 12
$ws.Cells.Item(12, 1).Value = 43053
$ws.Cells.Item(12, 2).Value = 2017
$ws.Cells.Item(12, 3).Value = 2.18621550610123
$ws.Cells.Item(12, 4).Value = 2018
$ws.Cells.Item(12, 5).Value = 2.066615940231964
# Row 13
$ws.Cells.Item(13, 1).Value = 43418
$ws.Cells.Item(13, 2).Value = 2018
$ws.Cells.Item(13, 3).Value = 1.911050033324102
$ws.Cells.Item(13, 4).Value = 2019
$ws.Cells.Item(13, 5).Value = 3.0862758122153
# Row 14
$ws.Cells.Item(14, 1).Value = 43783
$ws.Cells.Item(14, 2).Value = 2019
$ws.Cells.Item(14, 3).Value = 1.457852003181337
$ws.Cells.Item(14, 4).Value = 2020
$ws.Cells.Item(14, 5).Value = -1.135072001636328
# Row 15
$ws.Cells.Item(15, 1).Value = 44159
$ws.Cells.Item(15, 2).Value = 2020
$ws.Cells.Item(15, 3).Value = -3.258619210312885
$ws.Cells.Item(15, 4).Value = 2021
$ws.Cells.Item(15, 5).Value = -2.878617960200258
# Row 16
$ws.Cells.Item(16, 1).Value = 44525
$ws.Cells.Item(16, 2).Value = 2021
$ws.Cells.Item(16, 3).Value = 0.4255262881966981
$ws.Cells.Item(16, 4).Value = 2022
$ws.Cells.Item(16, 5).Value = 1.466936654457096
# Row 17
$ws.Cells.Item(17, 1).Value = 44890
$ws.Cells.Item(17, 2).Value = 2022
$ws.Cells.Item(17, 3).Value = 3.293290997728171
$ws.Cells.Item(17, 4).Value = 2023
$ws.Cells.Item(17, 5).Value = -1.421977974472588
# Row 18
$ws.Cells.Item(18, 1).Value = 45254
$ws.Cells.Item(18, 2).Value = 2023
$ws.Cells.Item(18, 3).Value = -0.2814561130375925
$ws.Cells.Item(18, 4).Value = 2024
$ws.Cells.Item(18, 5).Value = -0.3873858053678236
# Row 19
$ws.Cells.Item(19, 1).Value = 45618
$ws.Cells.Item(19, 2).Value = 2024
$ws.Cells.Item(19, 3).Value = -0.6470065423293869
$ws.Cells.Item(19, 4).Value = 2025
$ws.Cells.Item(19, 5).Value = 1.276847713071927
Write-Host "Forecast series refreshed: rows 2-19 updated"
